# [IMP] Adjust template CD Receivabl Planning
# 1. Make the "Receivable CD Planning Detail Report" sheet visible again
# 2. Fix "Quater" typo -> "Quarter" in the quarter headers on both sheets

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Receivable CD Planning Report")
$ws1.Range("B12").Value = "Quarter 1"
$ws1.Range("C12").Value = "Quarter 2"
$ws1.Range("D12").Value = "Quarter 3"
$ws1.Range("E12").Value = "Quarter 4"

$ws2 = $wb.Worksheets.Item("Receivable CD Planning Detail Report")
$ws2.Range("F12").Value = "Quarter 1"
$ws2.Range("G12").Value = "Quarter 2"
$ws2.Range("H12").Value = "Quarter 3"
$ws2.Range("I12").Value = "Quarter 4"

# Unhide the detail report sheet
$ws2.Visible = $true
